# Generate Report for Handoff
# Replace the old handoff id (c34dcc64-3751-406a-a388-432951ab0ff0) with the
# new one (ad23a82d-6c26-4571-906c-12a7624f43dc) everywhere it appears as a
# cell value / hyperlink display text, and bump the associated timestamps.

$wb = $excel.ActiveWorkbook

$oldId = "c34dcc64-3751-406a-a388-432951ab0ff0"
$newId = "ad23a82d-6c26-4571-906c-12a7624f43dc"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("C2").Value = "2016-50-21 00:50:53"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d34ec1d46d468b0c23af873b2b61d06ac4bb60a1/e2e/$oldId.md",
    "",
    "",
    "$newId.md"
)

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("D2").Value = "$newId.e82153bd6d0b6c34468dad8a28d71f51695be75f.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-21 00:50:49"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d34ec1d46d468b0c23af873b2b61d06ac4bb60a1/e2e/$oldId.md",
    "",
    "",
    "$newId.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d34ec1d46d468b0c23af873b2b61d06ac4bb60a1/e2e/$oldId.md",
    "",
    "",
    ".md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d3b0bc71595e6785b641b4d6f78ac1c255a7be8b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldId.4e53e07e232566d0a0b66d51a0352c65eb4c957b.zh-cn.xlf",
    "",
    "",
    "$newId.e82153bd6d0b6c34468dad8a28d71f51695be75f.zh-cn.xlf"
)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("D2").Value = "$newId.e82153bd6d0b6c34468dad8a28d71f51695be75f.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-21 00:50:53"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d34ec1d46d468b0c23af873b2b61d06ac4bb60a1/e2e/$oldId.md",
    "",
    "",
    "$newId.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d34ec1d46d468b0c23af873b2b61d06ac4bb60a1/e2e/$oldId.md",
    "",
    "",
    ".md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/af0a2f3a24427d5bb355c7f3ebe5464cbaf4b377/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldId.4e53e07e232566d0a0b66d51a0352c65eb4c957b.de-de.xlf",
    "",
    "",
    "$newId.e82153bd6d0b6c34468dad8a28d71f51695be75f.de-de.xlf"
)
